$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 11007.294
$ws.Range("J58").Value = 45008.5
$ws.Range("L58").Value = 135025.5
$ws.Range("N58").Value = -135325.5
$ws.Range("H82").Value = 16870.75
$ws.Range("I82").Value = 16870.75
$ws.Range("K82").Value = 50612.25
$ws.Range("M82").Value = -50206.25
$ws.Range("H85").Value = 16870.75
$ws.Range("I85").Value = 16870.75
$ws.Range("K85").Value = 50612.25
$ws.Range("M85").Value = -49208.25
$ws.Range("H92").Value = 961.5217
$ws.Range("I92").Value = 1074.0625
$ws.Range("J92").Value = 704.2857
$ws.Range("K92").Value = 1074.0625
$ws.Range("L92").Value = 704.2857
$ws.Range("M92").Value = 173.9375
$ws.Range("N92").Value = -3200.2857
$ws.Range("H94").Value = 877
$ws.Range("I94").Value = 877
$ws.Range("K94").Value = 877
$ws.Range("M94").Value = -426
$ws.Range("H116").Value = 74171.78999999999
$ws.Range("I116").Value = 113478.336
$ws.Range("K116").Value = 113478.336
$ws.Range("M116").Value = -110036.336
$ws.Range("H134").Value = 89999.2
$ws.Range("J134").Value = 89999.2
$ws.Range("L134").Value = 89999.2
$ws.Range("N134").Value = -100139.2
$ws.Range("H137").Value = 2002.25
$ws.Range("I137").Value = 1296.8462
$ws.Range("K137").Value = 3890.5386
$ws.Range("M137").Value = -1340.5386
$ws.Range("H141").Value = 75076.766
$ws.Range("I141").Value = 75076.766
$ws.Range("K141").Value = 225230.298
$ws.Range("M141").Value = -220050.298

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 27172.25
$ws.Range("J24").Value = 27172.25
$ws.Range("L24").Value = 27172.25
$ws.Range("N24").Value = -27920.25
$ws.Range("H43").Value = 30207.125
$ws.Range("J43").Value = 31194.334
$ws.Range("L43").Value = 31194.334
$ws.Range("N43").Value = -31820.334
$ws.Range("H63").Value = 13975.223
$ws.Range("I63").Value = 14791.25
$ws.Range("J63").Value = 7447
$ws.Range("K63").Value = 14791.25
$ws.Range("L63").Value = 7447
$ws.Range("M63").Value = -14105.25
$ws.Range("N63").Value = -8819
$ws.Range("H66").Value = 13975.223
$ws.Range("I66").Value = 14791.25
$ws.Range("J66").Value = 7447
$ws.Range("K66").Value = 73956.25
$ws.Range("L66").Value = 37235
$ws.Range("M66").Value = -70524.25
$ws.Range("N66").Value = -44099
$ws.Range("H74").Value = 3020.7026
$ws.Range("I74").Value = 3130.3333
$ws.Range("K74").Value = 3130.3333
$ws.Range("M74").Value = -2256.3333
$ws.Range("H77").Value = 3020.7026
$ws.Range("I77").Value = 3130.3333
$ws.Range("K77").Value = 15651.6665
$ws.Range("M77").Value = -11283.6665
$ws.Range("H88").Value = 2270.6875
$ws.Range("I88").Value = 2014.6
$ws.Range("J88").Value = 2387.0908
$ws.Range("K88").Value = 2014.6
$ws.Range("L88").Value = 2387.0908
$ws.Range("M88").Value = -1608.6
$ws.Range("N88").Value = -3199.0908
$ws.Range("H91").Value = 2270.6875
$ws.Range("I91").Value = 2014.6
$ws.Range("J91").Value = 2387.0908
$ws.Range("K91").Value = 2014.6
$ws.Range("L91").Value = 2387.0908
$ws.Range("M91").Value = -610.5999999999999
$ws.Range("N91").Value = -5195.0908
$ws.Range("H96").Value = 28276.666
$ws.Range("J96").Value = 28276.666
$ws.Range("L96").Value = 28276.666
$ws.Range("N96").Value = -33768.666
$ws.Range("H100").Value = 27172.25
$ws.Range("J100").Value = 27172.25
$ws.Range("L100").Value = 27172.25
$ws.Range("N100").Value = -29336.25
$ws.Range("H110").Value = 4703.385
$ws.Range("J110").Value = 7967.2856
$ws.Range("L110").Value = 7967.2856
$ws.Range("N110").Value = -12057.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 10824.36
$ws.Range("I94").Value = 3218.1177
$ws.Range("J94").Value = 26987.625
$ws.Range("K94").Value = 3218.1177
$ws.Range("L94").Value = 26987.625
$ws.Range("M94").Value = -2767.1177
$ws.Range("N94").Value = -27889.625
$ws.Range("H99").Value = 159543
$ws.Range("I99").Value = 5806.5
$ws.Range("K99").Value = 5806.5
$ws.Range("M99").Value = -4308.5
$ws.Range("H107").Value = 3188.3076
$ws.Range("I107").Value = 2403.182
$ws.Range("K107").Value = 2403.182
$ws.Range("M107").Value = -483.1819999999998
$ws.Range("H134").Value = 17680092
$ws.Range("I134").Value = 8336791
$ws.Range("K134").Value = 25010373
$ws.Range("M134").Value = -25007838

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2301.625
$ws.Range("I16").Value = 2078
$ws.Range("J16").Value = 2525.25
$ws.Range("K16").Value = 2078
$ws.Range("L16").Value = 2525.25
$ws.Range("M16").Value = -1791
$ws.Range("N16").Value = -3099.25
$ws.Range("H113").Value = 2301.625
$ws.Range("I113").Value = 2078
$ws.Range("J113").Value = 2525.25
$ws.Range("K113").Value = 2078
$ws.Range("L113").Value = 2525.25
$ws.Range("M113").Value = 92
$ws.Range("N113").Value = -6865.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 754.2258
$ws.Range("I113").Value = 501
$ws.Range("J113").Value = 771.6896400000001
$ws.Range("K113").Value = 1503
$ws.Range("L113").Value = 2315.06892
$ws.Range("M113").Value = 667
$ws.Range("N113").Value = -6655.06892
$ws.Range("H116").Value = 140001.5
$ws.Range("I116").Value = 197859.58
$ws.Range("J116").Value = 4999.3335
$ws.Range("K116").Value = 593578.74
$ws.Range("L116").Value = 14998.0005
$ws.Range("M116").Value = -590136.74
$ws.Range("N116").Value = -21882.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2331.3044
$ws.Range("J113").Value = 2172.7144
$ws.Range("L113").Value = 2172.7144
$ws.Range("N113").Value = -6512.7144
$ws.Range("H117").Value = 28999.5
$ws.Range("J117").Value = 28999.5
$ws.Range("L117").Value = 28999.5
$ws.Range("N117").Value = -35883.5
$ws.Range("H132").Value = 2029.2106
$ws.Range("I132").Value = 2058.9443
$ws.Range("J132").Value = 1494
$ws.Range("K132").Value = 6176.8329
$ws.Range("L132").Value = 4482
$ws.Range("M132").Value = -3646.8329
$ws.Range("N132").Value = -9542

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3534.2727
$ws.Range("I7").Value = 3208.5557
$ws.Range("K7").Value = 3208.5557
$ws.Range("M7").Value = -3096.5557
$ws.Range("H61").Value = 2252.111
$ws.Range("I61").Value = 1533
$ws.Range("K61").Value = 1533
$ws.Range("M61").Value = -1331
$ws.Range("H93").Value = 16719.207
$ws.Range("I93").Value = 2067.1875
$ws.Range("J93").Value = 34752.46
$ws.Range("K93").Value = 2067.1875
$ws.Range("L93").Value = 34752.46
$ws.Range("M93").Value = -819.1875
$ws.Range("N93").Value = -37248.46
$ws.Range("H113").Value = 2252.111
$ws.Range("I113").Value = 1533
$ws.Range("K113").Value = 1533
$ws.Range("M113").Value = 637
$ws.Range("H126").Value = 3534.2727
$ws.Range("I126").Value = 3208.5557
$ws.Range("K126").Value = 9625.667099999999
$ws.Range("M126").Value = -7155.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3439
$ws.Range("I81").Value = 1657.8334
$ws.Range("K81").Value = 3315.6668
$ws.Range("M81").Value = -2254.6668
$ws.Range("H84").Value = 3439
$ws.Range("I84").Value = 1657.8334
$ws.Range("K84").Value = 16578.334
$ws.Range("M84").Value = -11274.334
$ws.Range("H96").Value = 4193.1816
$ws.Range("J96").Value = 3692.4443
$ws.Range("L96").Value = 3692.4443
$ws.Range("N96").Value = -6438.4443
